# Apply the "Q5"/"W5" quiz columns (M, N) to the grades sheet.
# Mirrors the OOXML diff: new shared strings "Q5"/"W5" in M1/N1,
# per-student quiz/exam values in M2:N33, and refreshed dimension/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (column M = Q5, column N = W5)
$ws.Range("M1").Value = "Q5"
$ws.Range("N1").Value = "W5"

$mnData = @{
    2 = @(2, 35)
    3 = @(0, 0)
    4 = @(0, 35)
    5 = @(2, 33)
    6 = @(2, 30)
    7 = @(0, 30)
    8 = @(2, 30)
    9 = @(2, 35)
    10 = @(2, 32)
    11 = @(2, 32)
    12 = @(2, 40)
    13 = @(2, 33)
    14 = @(4, 40)
    15 = @(2, 40)
    16 = @(2, 30)
    17 = @(2, 32)
    18 = @(4, 35)
    19 = @(4, 33)
    20 = @(2, 35)
    21 = @(2, 33)
    22 = @(0, 30)
    23 = @(4, 30)
    24 = @(4, 35)
    25 = @(4, 30)
    26 = @(4, 30)
    27 = @(4, 35)
    28 = @(2, 40)
    29 = @(0, 40)
    30 = @(4, 35)
    31 = @(0, 33)
    32 = @(4, 33)
    33 = @(0, 0)
}

foreach ($r in $mnData.Keys) {
    $vals = $mnData[$r]
    $ws.Cells.Item($r, 13).Value = $vals[0]
    $ws.Cells.Item($r, 14).Value = $vals[1]
}

# Update the active selection to N33, matching the post-edit workbook state.
$ws.Range("N33").Select() | Out-Null
